$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "77÷4="
$t.Cell(1, 2).Range.Text = "84÷4="
$t.Cell(1, 3).Range.Text = "58÷5="
$t.Cell(1, 4).Range.Text = "93÷4="
$t.Cell(1, 5).Range.Text = "58÷2="

$t.Cell(5, 1).Range.Text = "35÷8="
$t.Cell(5, 2).Range.Text = "70÷3="
$t.Cell(5, 3).Range.Text = "70÷7="
$t.Cell(5, 4).Range.Text = "77÷5="
$t.Cell(5, 5).Range.Text = "14÷4="

$t.Cell(9, 1).Range.Text = "95÷8="
$t.Cell(9, 2).Range.Text = "61÷5="
$t.Cell(9, 3).Range.Text = "70÷3="
$t.Cell(9, 4).Range.Text = "85÷7="
$t.Cell(9, 5).Range.Text = "86÷6="

$t.Cell(13, 1).Range.Text = "98÷2="
$t.Cell(13, 2).Range.Text = "60÷7="
$t.Cell(13, 3).Range.Text = "73÷4="
$t.Cell(13, 4).Range.Text = "55÷6="
$t.Cell(13, 5).Range.Text = "22÷4="

$t.Cell(17, 1).Range.Text = "39÷4="
$t.Cell(17, 2).Range.Text = "21÷7="
$t.Cell(17, 3).Range.Text = "93÷2="
$t.Cell(17, 4).Range.Text = "76÷3="
$t.Cell(17, 5).Range.Text = "44÷8="
